$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.906.35'
$ws.Range('E2').Value = '  -2.19%  '
$ws.Range('D3').Value = '1.652.86'
$ws.Range('E3').Value = '  -0.97%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.41'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3891'
$ws.Range('E7').Value = '  -1.41%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3826'
$ws.Range('E8').Value = '  -2.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '51.78'
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.344'
$ws.Range('E10').Value = '  -3.29%  '
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08462'
$ws.Range('E12').Value = '  -1.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '23.96'
$ws.Range('E13').Value = '  -2.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.035'
$ws.Range('E14').Value = '  -3.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.030'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001315'
$ws.Range('E16').Value = '  -1.53%  '
$ws.Range('D17').Value = '1.655.01'
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '94.18'
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06990'
$ws.Range('E19').Value = '  -0.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.62'
$ws.Range('E20').Value = '  -4.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.971'
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('E23').Value = '  -0.38%  '
$ws.Range('D24').Value = '23.897.64'
$ws.Range('E24').Value = '  -2.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.449'
$ws.Range('E25').Value = '  -1.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.949'
$ws.Range('E26').Value = '  -4.87%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.03'
$ws.Range('E27').Value = '  -2.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '153.34'
$ws.Range('E28').Value = '  -2.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.417'
$ws.Range('E29').Value = '  +1.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '137.46'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.745'
$ws.Range('E31').Value = '  -2.47%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.482'
$ws.Range('E32').Value = '  -2.69%  '
$ws.Range('D33').Value = '1.834.81'
$ws.Range('E33').Value = '  -0.81%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08142'
$ws.Range('E34').Value = '  -1.28%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9992'
$ws.Range('E35').Value = '  -6.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02909'
$ws.Range('E36').Value = '  -6.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.668'
$ws.Range('E37').Value = '  -3.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '10.74'
$ws.Range('E38').Value = '  -3.64%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2675'
$ws.Range('E39').Value = '  -3.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.09111'
$ws.Range('E40').Value = '  -1.61%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.7573'
$ws.Range('E41').Value = '  -1.58%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '13.56'
$ws.Range('E42').Value = '  -1.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.421'
$ws.Range('E43').Value = '  -1.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.44'
$ws.Range('E44').Value = '  -0.98%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6938'
$ws.Range('E45').Value = '  -2.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.448'
$ws.Range('E46').Value = '  -3.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.113'
$ws.Range('E47').Value = '  -0.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9999'
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08276'
$ws.Range('E49').Value = '  -1.82%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '133.63'
$ws.Range('E50').Value = '  -2.23%  '
$ws.Range('E51').Value = '  -2.73%  '
